# Auto-generated edit script applying cryptos list update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '''69.321.05'
$ws.Range('E2').Value = '  +1.44%  '
$ws.Range('D3').Value = '''3.901.17'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.13%  '
$ws.Range('D5').Value = '''528.92'
$ws.Range('E5').Value = '  +8.75%  '
$ws.Range('D6').Value = '''144.89'
$ws.Range('E6').Value = '  -0.76%  '
$ws.Range('D7').Value = '''0.614'
$ws.Range('E7').Value = '  -1.44%  '
$ws.Range('E8').Value = '  +0.09%  '
$ws.Range('E9').Value = '  -2.96%  '
$ws.Range('D10').Value = '''0.172'
$ws.Range('E10').Value = '  -5.52%  '
$ws.Range('D11').Value = '''0.0000334'
$ws.Range('E11').Value = '  -5.88%  '
$ws.Range('D12').Value = '''42.08'
$ws.Range('E12').Value = '  -2.05%  '
$ws.Range('D13').Value = '''4.521.66'
$ws.Range('E13').Value = '  +0.01%  '
$ws.Range('E14').Value = '  -2.00%  '
$ws.Range('D15').Value = '''3.890.99'
$ws.Range('E15').Value = '  +0.06%  '
$ws.Range('D16').Value = '''14.02'
$ws.Range('E16').Value = '  -1.35%  '
$ws.Range('E17').Value = '  +6.70%  '
$ws.Range('E18').Value = '  -1.46%  '
$ws.Range('D19').Value = '''19.78'
$ws.Range('D20').Value = '''69.283.60'
$ws.Range('E20').Value = '  +1.40%  '
$ws.Range('D21').Value = '''424.73'
$ws.Range('E21').Value = '  -1.51%  '
$ws.Range('E22').Value = '  -5.46%  '
$ws.Range('D23').Value = '''14.20'
$ws.Range('E23').Value = '  -3.98%  '
$ws.Range('D24').Value = '''87.88'
$ws.Range('E24').Value = '  -1.54%  '
$ws.Range('D25').Value = '''4.01'
$ws.Range('E25').Value = '  +7.94%  '
$ws.Range('D26').Value = '''11.41'
$ws.Range('E26').Value = '  -8.68%  '
$ws.Range('D27').Value = '''10.63'
$ws.Range('E27').Value = '  -3.22%  '
$ws.Range('D28').Value = '''36.39'
$ws.Range('E28').Value = '  -2.48%  '
$ws.Range('D29').Value = '''689.34'
$ws.Range('E29').Value = '  -4.65%  '
$ws.Range('D30').Value = '''13.23'
$ws.Range('E30').Value = '  -1.52%  '
$ws.Range('E31').Value = '  -3.14%  '
$ws.Range('E32').Value = '  -2.58%  '
$ws.Range('D33').Value = '''68.29'
$ws.Range('E33').Value = '  +10.71%  '
$ws.Range('E34').Value = '  +8.32%  '
$ws.Range('D35').Value = '''0.0₃0861'
$ws.Range('E35').Value = '  -3.11%  '
$ws.Range('D36').Value = '''5.91'
$ws.Range('E36').Value = '  -2.67%  '
$ws.Range('D37').Value = '''40.03'
$ws.Range('E37').Value = '  -1.97%  '
$ws.Range('D38').Value = '''0.149'
$ws.Range('E38').Value = '  +1.31%  '
$ws.Range('E39').Value = '  +0.14%  '
$ws.Range('E40').Value = '  +8.74%  '
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('E42').Value = '  -2.14%  '
$ws.Range('E43').Value = '  +7.08%  '
$ws.Range('E44').Value = '  -6.79%  '
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('E46').Value = '  -1.22%  '
$ws.Range('D47').Value = '''0.000279'
$ws.Range('E47').Value = '  +13.11%  '
$ws.Range('E48').Value = '  +6.68%  '
$ws.Range('D49').Value = '''2.749.09'
$ws.Range('E49').Value = '  +14.80%  '
$ws.Range('B50').Value = 'Monero'
$ws.Range('C50').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D50').Value = '''144.59'
$ws.Range('E50').Value = '  +0.32%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = '''26.48'
$ws.Range('E51').Value = '  +5.82%  '
